$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.546.14'
$ws.Range('E2').Value = '  +0.53%  '
$ws.Range('D3').Value = '2.136.80'
$ws.Range('E3').Value = '  +1.77%  '
$ws.Range('D4').Value = "'1.009"
$ws.Range('E4').Value = '  +0.55%  '
$ws.Range('D5').Value = "'352.22"
$ws.Range('E5').Value = '  +5.45%  '
$ws.Range('D6').Value = "'1.008"
$ws.Range('E6').Value = '  +0.50%  '
$ws.Range('D7').Value = "'0.5258"
$ws.Range('E7').Value = '  +1.07%  '
$ws.Range('D8').Value = "'0.4552"
$ws.Range('E8').Value = '  +0.45%  '
$ws.Range('D9').Value = "'53.58"
$ws.Range('E9').Value = '  -1.51%  '
$ws.Range('D10').Value = "'0.09182"
$ws.Range('E10').Value = '  +3.50%  '
$ws.Range('D11').Value = "'1.184"
$ws.Range('E11').Value = '  +0.43%  '
$ws.Range('D12').Value = "'25.33"
$ws.Range('E12').Value = '  +5.40%  '
$ws.Range('D13').Value = '2.140.46'
$ws.Range('E13').Value = '  +2.25%  '
$ws.Range('D14').Value = "'6.903"
$ws.Range('E14').Value = '  +1.68%  '
$ws.Range('D15').Value = "'8.187"
$ws.Range('E15').Value = '  +2.31%  '
$ws.Range('D16').Value = "'102.27"
$ws.Range('E16').Value = '  +5.75%  '
$ws.Range('E17').Value = '  +2.87%  '
$ws.Range('D18').Value = "'1.009"
$ws.Range('E18').Value = '  +0.38%  '
$ws.Range('D19').Value = "'0.06716"
$ws.Range('E19').Value = '  +1.43%  '
$ws.Range('D20').Value = "'20.34"
$ws.Range('E20').Value = '  +6.23%  '
$ws.Range('D21').Value = "'1.007"
$ws.Range('E21').Value = '  +0.43%  '
$ws.Range('D22').Value = "'6.376"
$ws.Range('E22').Value = '  +1.80%  '
$ws.Range('D23').Value = '30.643.25'
$ws.Range('E23').Value = '  +0.67%  '
$ws.Range('D24').Value = "'12.89"
$ws.Range('E24').Value = '  +4.74%  '
$ws.Range('D25').Value = "'2.385"
$ws.Range('E25').Value = '  +2.22%  '
$ws.Range('D26').Value = '2.391.35'
$ws.Range('E26').Value = '  +2.24%  '
$ws.Range('D27').Value = "'2.658"
$ws.Range('E27').Value = '  +6.05%  '
$ws.Range('D28').Value = "'22.48"
$ws.Range('E28').Value = '  +1.61%  '
$ws.Range('D29').Value = "'164.72"
$ws.Range('E29').Value = '  +1.38%  '
$ws.Range('D30').Value = "'136.84"
$ws.Range('E30').Value = '  +3.04%  '
$ws.Range('D31').Value = "'1.225"
$ws.Range('E31').Value = '  +2.04%  '
$ws.Range('D32').Value = "'0.1083"
$ws.Range('E32').Value = '  +1.62%  '
$ws.Range('D33').Value = "'1.685"
$ws.Range('E33').Value = '  +2.09%  '
$ws.Range('D34').Value = "'6.401"
$ws.Range('E34').Value = '  +0.69%  '
$ws.Range('D35').Value = "'4.042"
$ws.Range('E35').Value = '  +2.48%  '
$ws.Range('D36').Value = "'6.158"
$ws.Range('E36').Value = '  +6.11%  '
$ws.Range('D37').Value = "'10.47"
$ws.Range('E37').Value = '  +1.21%  '
$ws.Range('E38').Value = '  +3.07%  '
$ws.Range('D39').Value = "'0.06974"
$ws.Range('E39').Value = '  +2.17%  '
$ws.Range('D40').Value = "'0.2335"
$ws.Range('E40').Value = '  +2.00%  '
$ws.Range('D41').Value = "'12.77"
$ws.Range('E41').Value = '  +0.45%  '
$ws.Range('D42').Value = "'0.6986"
$ws.Range('E42').Value = '  +1.93%  '
$ws.Range('D43').Value = "'1.283"
$ws.Range('E43').Value = '  +3.26%  '
$ws.Range('D44').Value = "'14.75"
$ws.Range('E44').Value = '  +6.03%  '
$ws.Range('D45').Value = "'2.354"
$ws.Range('E45').Value = '  +1.49%  '
$ws.Range('D46').Value = "'0.6503"
$ws.Range('E46').Value = '  +2.64%  '
$ws.Range('D47').Value = "'0.00000000371"
$ws.Range('E47').Value = '  +8.08%  '
$ws.Range('D48').Value = "'3.752"
$ws.Range('E48').Value = '  +2.56%  '
$ws.Range('D49').Value = "'1.253"
$ws.Range('E49').Value = '  +0.66%  '
$ws.Range('D50').Value = "'84.11"
$ws.Range('E50').Value = '  +1.68%  '
$ws.Range('D51').Value = "'0.07297"
$ws.Range('E51').Value = '  +2.51%  '
